# "Generate Report for Archive"
# The localization status moved on from "Ready for handoff" -> "In Translation"
# for this entry. Update every sheet that surfaces the status for this file
# (the Overview rollup columns per-locale, and each locale's own Status
# column), then re-fit the affected columns now that the text is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview: the zh-cn / de-de columns hold each locale's current status.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-locale sheets: Status column (C).
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# The status text got shorter, so the columns that show it shrink to fit.
$fitWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $fitWidth
$wsOverview.Columns.Item(6).ColumnWidth = $fitWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $fitWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $fitWidth
